$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Save off the current (pre-edit) values we need to move around ---
$groupsHeader = $ws.Range("E1").Value2     # "GROUPS"
$keyHeader    = $ws.Range("H1").Value2     # "KEY"

$groupsRow2 = $ws.Range("E2").Value2       # "production, test"
$keyRow3    = $ws.Range("H3").Value2       # SSH public key text
$groupsRow4 = $ws.Range("E4").Value2       # "preproduction"

# --- Header row: KEY moves to column E, GROUPS moves to column H, new COMPLIANCE_GROUPS in I ---
$ws.Range("E1").Value2 = $keyHeader
$ws.Range("H1").Value2 = $groupsHeader
$ws.Range("I1").Value2 = "COMPLIANCE_GROUPS"

# --- Row 2 ---
# E2 (old GROUPS cell) becomes an empty but still-formatted KEY cell, so give it
# the same look the rest of the credential cells on this row carry.
$ws.Range("A2").Style = $ws.Range("F2").Style
$ws.Range("E2").ClearContents()
$ws.Range("E2").Style = $ws.Range("F2").Style
$ws.Range("H2").Value2 = $groupsRow2
$ws.Range("I2").Value2 = "anssi"

# --- Row 3 ---
$ws.Range("E3").Value2 = $keyRow3
$ws.Range("E3").WrapText = $true
$ws.Range("F3").Value2 = "user"
$ws.Range("H3").Clear()
$ws.Range("I3").Value2 = "anssi"

# --- Row 4 ---
$ws.Range("E4").Clear()
$ws.Range("H4").Value2 = $groupsRow4
$ws.Range("I4").Value2 = "anssi"

# --- Column widths (characters) ---
$ws.Columns.Item(3).ColumnWidth = 16.65
$ws.Columns.Item(5).ColumnWidth = 7.16
$ws.Columns.Item(6).ColumnWidth = 12.68
$ws.Columns.Item(7).ColumnWidth = 15.32
$ws.Columns.Item(8).ColumnWidth = 14.44
$ws.Columns.Item(9).ColumnWidth = 20.06

# --- Row heights ---
$ws.Rows.Item(3).RowHeight = 163.5

Write-Output "done"
